# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the per-job worksheets with freshly
# pulled prices, recomputed leve costs and profits.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1321.1666
$ws.Range("I28").Value = 565.2857
$ws.Range("K28").Value = 565.2857
$ws.Range("M28").Value = -80.28570000000002

$ws.Range("H32").Value = 13889977
$ws.Range("I32").Value = 27778378
$ws.Range("J32").Value = 1575.6666
$ws.Range("K32").Value = 27778378
$ws.Range("L32").Value = 1575.6666
$ws.Range("M32").Value = -27778052
$ws.Range("N32").Value = -2227.6666

$ws.Range("H43").Value = 2349.5
$ws.Range("J43").Value = 2349.5
$ws.Range("L43").Value = 2349.5
$ws.Range("N43").Value = -2487.5

$ws.Range("H137").Value = 2692.577
$ws.Range("I137").Value = 1581.0555
$ws.Range("J137").Value = 5193.5
$ws.Range("K137").Value = 4743.166499999999
$ws.Range("L137").Value = 15580.5
$ws.Range("M137").Value = -2193.166499999999
$ws.Range("N137").Value = -20680.5

$ws.Range("H138").Value = 2432.0278
$ws.Range("I138").Value = 1548.1428
$ws.Range("J138").Value = 2645.3794
$ws.Range("K138").Value = 4644.428400000001
$ws.Range("L138").Value = 7936.138199999999
$ws.Range("M138").Value = 495.5715999999993
$ws.Range("N138").Value = -18216.1382

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2230027.8
$ws.Range("I32").Value = 1124803.4
$ws.Range("K32").Value = 1124803.4
$ws.Range("M32").Value = -1124516.4

$ws.Range("H61").Value = 2857.8708
$ws.Range("I61").Value = 2611.45
$ws.Range("K61").Value = 2611.45
$ws.Range("M61").Value = -2399.45

$ws.Range("H74").Value = 2245.2917
$ws.Range("I74").Value = 1893.8889
$ws.Range("J74").Value = 3299.5
$ws.Range("K74").Value = 1893.8889
$ws.Range("L74").Value = 3299.5
$ws.Range("M74").Value = -1019.8889
$ws.Range("N74").Value = -5047.5

$ws.Range("H77").Value = 2245.2917
$ws.Range("I77").Value = 1893.8889
$ws.Range("J77").Value = 3299.5
$ws.Range("K77").Value = 9469.4445
$ws.Range("L77").Value = 16497.5
$ws.Range("M77").Value = -5101.4445
$ws.Range("N77").Value = -25233.5

$ws.Range("H102").Value = 2424.389
$ws.Range("I102").Value = 1914.1
$ws.Range("K102").Value = 1914.1
$ws.Range("M102").Value = -292.0999999999999

$ws.Range("H132").Value = 12874.75
$ws.Range("I132").Value = 8599.6
$ws.Range("K132").Value = 25798.8
$ws.Range("M132").Value = -23268.8

$ws.Range("H136").Value = 2857.8708
$ws.Range("I136").Value = 2611.45
$ws.Range("K136").Value = 7834.349999999999
$ws.Range("M136").Value = -5284.349999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 85877.336
$ws.Range("I99").Value = 102102.9
$ws.Range("J99").Value = 4749.5
$ws.Range("K99").Value = 102102.9
$ws.Range("L99").Value = 4749.5
$ws.Range("M99").Value = -100604.9
$ws.Range("N99").Value = -7745.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 34.333332
$ws.Range("I35").Value = 25
$ws.Range("J35").Value = 39
$ws.Range("K35").Value = 25
$ws.Range("L35").Value = 39
$ws.Range("M35").Value = 269
$ws.Range("N35").Value = -627

$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 5000
$ws.Range("K45").Value = 5000
$ws.Range("M45").Value = -4407

$ws.Range("H80").Value = 49999
$ws.Range("J80").Value = 49999
$ws.Range("L80").Value = 49999
$ws.Range("N80").Value = -52245

$ws.Range("H83").Value = 49999
$ws.Range("J83").Value = 49999
$ws.Range("L83").Value = 149997
$ws.Range("N83").Value = -161229

$ws.Range("H99").Value = 3507.8667
$ws.Range("I99").Value = 3093.3333
$ws.Range("J99").Value = 4129.6665
$ws.Range("K99").Value = 3093.3333
$ws.Range("L99").Value = 4129.6665
$ws.Range("M99").Value = -1595.3333
$ws.Range("N99").Value = -7125.6665

$ws.Range("H126").Value = 3507.8667
$ws.Range("I126").Value = 3093.3333
$ws.Range("J126").Value = 4129.6665
$ws.Range("K126").Value = 9279.999899999999
$ws.Range("L126").Value = 12388.9995
$ws.Range("M126").Value = -6809.999899999999
$ws.Range("N126").Value = -17328.9995

$ws.Range("H132").Value = 9528360
$ws.Range("I132").Value = 3213.1428
$ws.Range("K132").Value = 9639.428400000001
$ws.Range("M132").Value = -7109.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 327.5
$ws.Range("I12").Value = 7
$ws.Range("J12").Value = 519.8
$ws.Range("K12").Value = 21
$ws.Range("L12").Value = 1559.4
$ws.Range("M12").Value = 152
$ws.Range("N12").Value = -1905.4

$ws.Range("H74").Value = 40006.5
$ws.Range("J74").Value = 40000
$ws.Range("L74").Value = 120000
$ws.Range("N74").Value = -122122

$ws.Range("H77").Value = 40006.5
$ws.Range("J77").Value = 40000
$ws.Range("L77").Value = 360000
$ws.Range("N77").Value = -370608

$ws.Range("H131").Value = 6949291.5
$ws.Range("J131").Value = 1982.8387
$ws.Range("L131").Value = 5948.5161
$ws.Range("N131").Value = -16028.5161

$ws.Range("H141").Value = 9303
$ws.Range("I141").Value = 3721.7273
$ws.Range("K141").Value = 11165.1819
$ws.Range("M141").Value = -5985.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1276.1666
$ws.Range("I97").Value = 1276.1666
$ws.Range("K97").Value = 1276.1666
$ws.Range("M97").Value = -780.1666

$ws.Range("H126").Value = 5889.75
$ws.Range("J126").Value = 11023.3
$ws.Range("L126").Value = 33069.89999999999
$ws.Range("N126").Value = -38009.89999999999

$ws.Range("H132").Value = 2456.4634
$ws.Range("I132").Value = 1846.762
$ws.Range("J132").Value = 3096.65
$ws.Range("K132").Value = 5540.286
$ws.Range("L132").Value = 9289.950000000001
$ws.Range("M132").Value = -3010.286
$ws.Range("N132").Value = -14349.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 34015.832
$ws.Range("I40").Value = 34015.832
$ws.Range("K40").Value = 34015.832
$ws.Range("M40").Value = -33879.832

$ws.Range("H46").Value = 3003.5454
$ws.Range("I46").Value = 2004.4445
$ws.Range("K46").Value = 2004.4445
$ws.Range("M46").Value = -1816.4445

$ws.Range("H80").Value = 49999
$ws.Range("J80").Value = 49999
$ws.Range("L80").Value = 49999
$ws.Range("N80").Value = -52245

$ws.Range("H83").Value = 49999
$ws.Range("J83").Value = 49999
$ws.Range("L83").Value = 149997
$ws.Range("N83").Value = -161229

$ws.Range("H136").Value = 2114.6191
$ws.Range("I136").Value = 1817.5
$ws.Range("K136").Value = 5452.5
$ws.Range("M136").Value = -2902.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4821.357
$ws.Range("I62").Value = 4374.75
$ws.Range("K62").Value = 4374.75
$ws.Range("M62").Value = -3750.75

$ws.Range("H65").Value = 4821.357
$ws.Range("I65").Value = 4374.75
$ws.Range("K65").Value = 21873.75
$ws.Range("M65").Value = -18753.75

$ws.Range("H113").Value = 404.08334
$ws.Range("I113").Value = 404.08334
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1212.25002
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 957.7499800000001
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 1967.8077
$ws.Range("I132").Value = 1584.7273
$ws.Range("K132").Value = 4754.1819
$ws.Range("M132").Value = -2224.1819
